$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 4.915454186040008
$ws.Cells.Item(2, 4).Value = 8.059116697391062
$ws.Cells.Item(2, 5).Value = 14.32885920019903
$ws.Cells.Item(2, 6).Value = 34.48813553848034
$ws.Cells.Item(2, 7).Value = 43.26678856483562
$ws.Cells.Item(2, 8).Value = 17.45222829421117
$ws.Cells.Item(2, 10).Value = 9.805364248503889
$ws.Cells.Item(2, 13).Value = 59.93228829467144

$ws.Cells.Item(3, 3).Value = 4.751297677005716
$ws.Cells.Item(3, 4).Value = 7.830008958316026
$ws.Cells.Item(3, 5).Value = 13.85707858474041
$ws.Cells.Item(3, 6).Value = 35.15516020029545
$ws.Cells.Item(3, 7).Value = 43.69454217062832
$ws.Cells.Item(3, 8).Value = 17.68170204920864
$ws.Cells.Item(3, 10).Value = 9.688160437001489
$ws.Cells.Item(3, 13).Value = 56.62984748817227

$ws.Cells.Item(4, 3).Value = 4.649697292912531
$ws.Cells.Item(4, 4).Value = 7.685580534485968
$ws.Cells.Item(4, 5).Value = 13.56249429853682
$ws.Cells.Item(4, 6).Value = 35.59566757088005
$ws.Cells.Item(4, 7).Value = 44.01337540754488
$ws.Cells.Item(4, 8).Value = 17.83380954292453
$ws.Cells.Item(4, 10).Value = 9.619371387037557
$ws.Cells.Item(4, 13).Value = 54.49656666426139

$ws.Cells.Item(5, 3).Value = 4.608160379575963
$ws.Cells.Item(5, 4).Value = 7.625843500384537
$ws.Cells.Item(5, 5).Value = 13.44138592253792
$ws.Cells.Item(5, 6).Value = 35.7827289625412
$ws.Cells.Item(5, 7).Value = 44.15686426090563
$ws.Cells.Item(5, 8).Value = 17.89855062843411
$ws.Cells.Item(5, 10).Value = 9.592156398688672
$ws.Cells.Item(5, 13).Value = 53.6010805822708

$ws.Cells.Item(6, 3).Value = 4.601257228816817
$ws.Cells.Item(6, 4).Value = 7.615872926794919
$ws.Cells.Item(6, 5).Value = 13.42121708392282
$ws.Cells.Item(6, 6).Value = 35.81423997601489
$ws.Cells.Item(6, 7).Value = 44.1814930354403
$ws.Cells.Item(6, 8).Value = 17.90946552452302
$ws.Cells.Item(6, 10).Value = 9.587687249535147
$ws.Cells.Item(6, 13).Value = 53.45081942142449

$ws.Cells.Item(7, 3).Value = 4.649137559061122
$ws.Cells.Item(7, 4).Value = 7.684778381178431
$ws.Cells.Item(7, 5).Value = 13.56086505434898
$ws.Cells.Item(7, 6).Value = 35.59816006603956
$ws.Cells.Item(7, 7).Value = 44.01525638528129
$ws.Cells.Item(7, 8).Value = 17.83467158509362
$ws.Cells.Item(7, 10).Value = 9.619001023692691
$ws.Cells.Item(7, 13).Value = 54.48459514400166

$ws.Cells.Item(8, 3).Value = 4.859060652288966
$ws.Cells.Item(8, 4).Value = 7.980935990993852
$ws.Cells.Item(8, 5).Value = 14.16730627235129
$ws.Cells.Item(8, 6).Value = 34.7115566877203
$ws.Cells.Item(8, 7).Value = 43.40227351399714
$ws.Cells.Item(8, 8).Value = 17.52898551761527
$ws.Cells.Item(8, 10).Value = 9.764305109521224
$ws.Cells.Item(8, 13).Value = 58.81572028645622

$ws.Cells.Item(9, 3).Value = 5.261526201115447
$ws.Cells.Item(9, 4).Value = 8.529398650333931
$ws.Cells.Item(9, 5).Value = 15.31084515464015
$ws.Cells.Item(9, 6).Value = 33.22985651592481
$ws.Cells.Item(9, 7).Value = 42.67136343001286
$ws.Cells.Item(9, 8).Value = 17.02146401379152
$ws.Cells.Item(9, 10).Value = 10.07356865933538
$ws.Cells.Item(9, 13).Value = 66.45909079625135

$ws.Cells.Item(10, 3).Value = 5.548263803092957
$ws.Cells.Item(10, 4).Value = 8.909675955913865
$ws.Cells.Item(10, 5).Value = 16.11525012659878
$ws.Cells.Item(10, 6).Value = 32.3150119436403
$ws.Cells.Item(10, 7).Value = 42.45552494866779
$ws.Cells.Item(10, 8).Value = 16.70900806236766
$ws.Cells.Item(10, 10).Value = 10.31442487838811
$ws.Cells.Item(10, 13).Value = 71.54554932629451

$ws.Cells.Item(11, 3).Value = 5.676166067931619
$ws.Cells.Item(11, 4).Value = 9.077218836521554
$ws.Cells.Item(11, 5).Value = 16.47207299712639
$ws.Cells.Item(11, 6).Value = 31.94075446909027
$ws.Cells.Item(11, 7).Value = 42.43399350275977
$ws.Cells.Item(11, 8).Value = 16.58101830945898
$ws.Cells.Item(11, 10).Value = 10.42670329079487
$ws.Cells.Item(11, 13).Value = 73.74354812867794

$ws.Cells.Item(12, 3).Value = 5.724192095648433
$ws.Cells.Item(12, 4).Value = 9.139841647955569
$ws.Cells.Item(12, 5).Value = 16.60578948980098
$ws.Cells.Item(12, 6).Value = 31.80544483277645
$ws.Cells.Item(12, 7).Value = 42.437394113351
$ws.Cells.Item(12, 8).Value = 16.53467891664123
$ws.Cells.Item(12, 10).Value = 10.46958885749049
$ws.Cells.Item(12, 13).Value = 74.55918116567652

$ws.Cells.Item(13, 3).Value = 5.713867598427241
$ws.Cells.Item(13, 4).Value = 9.126391834969201
$ws.Cells.Item(13, 5).Value = 16.57705509129294
$ws.Cells.Item(13, 6).Value = 31.83429446403534
$ws.Cells.Item(13, 7).Value = 42.43613979111277
$ws.Cells.Item(13, 8).Value = 16.54456277184991
$ws.Cells.Item(13, 10).Value = 10.46033664790615
$ws.Cells.Item(13, 13).Value = 74.38426305060062

$ws.Cells.Item(14, 3).Value = 5.680125598990881
$ws.Cells.Item(14, 4).Value = 9.082387533851865
$ws.Cells.Item(14, 5).Value = 16.48310258039757
$ws.Cells.Item(14, 6).Value = 31.92949160616219
$ws.Cells.Item(14, 7).Value = 42.43403879854488
$ws.Cells.Item(14, 8).Value = 16.57716274862739
$ws.Cells.Item(14, 10).Value = 10.43022425577496
$ws.Cells.Item(14, 13).Value = 73.81098557082314

$ws.Cells.Item(15, 3).Value = 5.659403354750148
$ws.Cells.Item(15, 4).Value = 9.055325501447815
$ws.Cells.Item(15, 5).Value = 16.42536844327523
$ws.Cells.Item(15, 6).Value = 31.98864982353337
$ws.Cells.Item(15, 7).Value = 42.43427171555527
$ws.Cells.Item(15, 8).Value = 16.59741109824921
$ws.Cells.Item(15, 10).Value = 10.41182685825058
$ws.Cells.Item(15, 13).Value = 73.45766097606068

$ws.Cells.Item(16, 3).Value = 5.53985025256956
$ws.Cells.Item(16, 4).Value = 8.898613608518685
$ws.Cells.Item(16, 5).Value = 16.09173923384965
$ws.Cells.Item(16, 6).Value = 32.34034629347623
$ws.Cells.Item(16, 7).Value = 42.45852077939599
$ws.Cells.Item(16, 8).Value = 16.71766502420883
$ws.Cells.Item(16, 10).Value = 10.30713991045031
$ws.Cells.Item(16, 13).Value = 71.39956917002849

$ws.Cells.Item(17, 3).Value = 5.465826402599496
$ws.Cells.Item(17, 4).Value = 8.801051778117873
$ws.Cells.Item(17, 5).Value = 15.88466427569868
$ws.Cells.Item(17, 6).Value = 32.56710089048758
$ws.Cells.Item(17, 7).Value = 42.49341869761574
$ws.Cells.Item(17, 8).Value = 16.7951251188746
$ws.Cells.Item(17, 10).Value = 10.2435967044659
$ws.Cells.Item(17, 13).Value = 70.10728071315653

$ws.Cells.Item(18, 3).Value = 5.423013366718439
$ws.Cells.Item(18, 4).Value = 8.744426100525789
$ws.Cells.Item(18, 5).Value = 15.76470801895726
$ws.Cells.Item(18, 6).Value = 32.70144395158615
$ws.Cells.Item(18, 7).Value = 42.52067629200982
$ws.Cells.Item(18, 8).Value = 16.84100390379932
$ws.Cells.Item(18, 10).Value = 10.20730484490095
$ws.Cells.Item(18, 13).Value = 69.35308531700791

$ws.Cells.Item(19, 3).Value = 5.408478392862195
$ws.Cells.Item(19, 4).Value = 8.725167161756268
$ws.Cells.Item(19, 5).Value = 15.72394965425791
$ws.Cells.Item(19, 6).Value = 32.74759212360743
$ws.Cells.Item(19, 7).Value = 42.5311225627466
$ws.Cells.Item(19, 8).Value = 16.85676262747802
$ws.Cells.Item(19, 10).Value = 10.19506177483642
$ws.Cells.Item(19, 13).Value = 69.0958566359071

$ws.Cells.Item(20, 3).Value = 5.473731170901959
$ws.Cells.Item(20, 4).Value = 8.811490551267124
$ws.Cells.Item(20, 5).Value = 15.90679669409255
$ws.Cells.Item(20, 6).Value = 32.5425540810798
$ws.Cells.Item(20, 7).Value = 42.48895614949463
$ws.Cells.Item(20, 8).Value = 16.78674148825014
$ws.Cells.Item(20, 10).Value = 10.25033459907246
$ws.Cells.Item(20, 13).Value = 70.24597610094273

$ws.Cells.Item(21, 3).Value = 5.690047826929261
$ws.Cells.Item(21, 4).Value = 9.095335260327341
$ws.Cells.Item(21, 5).Value = 16.51073752158328
$ws.Cells.Item(21, 6).Value = 31.90135260565309
$ws.Cells.Item(21, 7).Value = 42.43433828241599
$ws.Cells.Item(21, 8).Value = 16.56752883578516
$ws.Cells.Item(21, 10).Value = 10.43905917319047
$ws.Cells.Item(21, 13).Value = 73.97982447469691

$ws.Cells.Item(22, 3).Value = 5.829030026992053
$ws.Cells.Item(22, 4).Value = 9.276037833264242
$ws.Cells.Item(22, 5).Value = 16.89722610690288
$ws.Cells.Item(22, 6).Value = 31.51987403292587
$ws.Cells.Item(22, 7).Value = 42.46620753526396
$ws.Cells.Item(22, 8).Value = 16.43670867104714
$ws.Cells.Item(22, 10).Value = 10.56453803426754
$ws.Cells.Item(22, 13).Value = 76.32276827879785

$ws.Cells.Item(23, 3).Value = 5.75508449401391
$ws.Cells.Item(23, 4).Value = 9.180044906103046
$ws.Cells.Item(23, 5).Value = 16.69173012342482
$ws.Cells.Item(23, 6).Value = 31.71990190339464
$ws.Cells.Item(23, 7).Value = 42.44285036775756
$ws.Cells.Item(23, 8).Value = 16.50535872572639
$ws.Cells.Item(23, 10).Value = 10.49737911292911
$ws.Cells.Item(23, 13).Value = 75.08120590209302

$ws.Cells.Item(24, 3).Value = 5.470158218393602
$ws.Cells.Item(24, 4).Value = 8.806772848712956
$ws.Cells.Item(24, 5).Value = 15.89679344266914
$ws.Cells.Item(24, 6).Value = 32.55363932939927
$ws.Cells.Item(24, 7).Value = 42.49095129956406
$ws.Cells.Item(24, 8).Value = 16.7905275423198
$ws.Cells.Item(24, 10).Value = 10.24728764778754
$ws.Cells.Item(24, 13).Value = 70.18330692748897

$ws.Cells.Item(25, 3).Value = 5.153999764704553
$ws.Cells.Item(25, 4).Value = 8.384811062073485
$ws.Cells.Item(25, 5).Value = 15.00725493998042
$ws.Cells.Item(25, 6).Value = 33.60151206343575
$ws.Cells.Item(25, 7).Value = 42.81510595569799
$ws.Cells.Item(25, 8).Value = 17.14848527227581
$ws.Cells.Item(25, 10).Value = 9.987421922039857
$ws.Cells.Item(25, 13).Value = 64.48373792389542
